# CHARMS - Verifying that audit trail history is visible to users + updates to test cases.
# Add "Name" rows alongside each Parent/Guardian/LAR "Signed" row on the
# rasE-ConsentAdult worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rasE-ConsentAdult")

# Insert top-down; row numbers below account for rows already inserted above.

# "Parent/Guardian 1 Name" goes right after "Parent/Guardian 1 Signed" (row 13)
$ws.Rows.Item(14).Insert()
$ws.Cells.Item(14, 1).Value = "Parent/Guardian 1 Name"
$ws.Rows.Item(14).RowHeight = 17

# "Parent/Guardian 2 Name" goes right after "Parent/Guardian 2 Signed" (now row 15)
$ws.Rows.Item(16).Insert()
$ws.Cells.Item(16, 1).Value = "Parent/Guardian 2 Name"
$ws.Rows.Item(16).RowHeight = 17

# "LAR 1 Name" goes right after "LAR 1 Signed" (now row 18)
$ws.Rows.Item(19).Insert()
$ws.Cells.Item(19, 1).Value = "LAR 1 Name"
$ws.Rows.Item(19).RowHeight = 17

# "LAR 2 Name" goes right after "LAR 2 Signed" (now row 20)
$ws.Rows.Item(21).Insert()
$ws.Cells.Item(21, 1).Value = "LAR 2 Name"
$ws.Rows.Item(21).RowHeight = 17
